$d = $word.ActiveDocument

# wdParagraph expand-unit constant.
$wdParagraph = 4

# Locate the paragraph that holds data row "4" (its leading spaces survive,
# the row content itself collapses to an ellipsis) and the paragraph that
# holds data row "98" (the first row to survive unchanged after the long
# run of deleted rows 5-97).
$row4Find = $d.Content
$row4Find.Find.ClearFormatting()
[void]$row4Find.Find.Execute("4      P     24.7    24.5    22.0    22.5", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $row4Find.Find.Found) {
    throw "could not find row 4"
}
[void]$row4Find.Expand($wdParagraph)
$row4Start = $row4Find.Start
$row4End = $row4Find.End

$row98Find = $d.Content
$row98Find.Find.ClearFormatting()
[void]$row98Find.Find.Execute("98      A     29.4    22.1    25.3     4.1", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $row98Find.Find.Found) {
    throw "could not find row 98"
}
[void]$row98Find.Expand($wdParagraph)
$row98Start = $row98Find.Start

# Delete everything between the end of row 4's paragraph and the start of
# row 98's paragraph - i.e. the paragraphs for rows 5 through 97, paragraph
# marks included.
$gap = $d.Range($row4End, $row98Start)
$gap.Delete()

# Row 4's paragraph now sits directly before row 98's. Collapse its visible
# data (everything after the 16-space indent) down to a single ellipsis
# character, leaving the leading whitespace run untouched.
$row4Para = $d.Range($row4Start, $row4Start)
[void]$row4Para.Expand($wdParagraph)
$row4Text = $row4Para.Text
$indent = $row4Text.Length - $row4Text.TrimStart(" ").Length
$dataRange = $d.Range($row4Para.Start + $indent, $row4Para.End - 1)
$dataRange.Text = "…"
